# edit.ps1
# Applies the weekly CompStat (79th Precinct) data refresh:
#   - header text: report period + volume/number bump
#   - per-crime-category Week/28-Day/YTD/2-Year figures and % changes
#   - a few rows (Transit/Shooting Vic./Shooting Inc./Hate Crimes) flip
#     cells between the "no data" placeholder text ("0" / "***.*") and
#     real numeric entries as reporting activity appears/disappears.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# 1. Header text: "Volume 32  Number 30" -> "...Number 31", and the
#    reporting week date range rolls forward one week.
# ---------------------------------------------------------------------
$ws.Range("A8").Value = "Volume 32   Number  31"
$ws.Range("C9").Value = "Report Covering the Week  7/28/2025  Through  8/3/2025"

# ---------------------------------------------------------------------
# 2. Plain numeric value updates (no type/style change).
# ---------------------------------------------------------------------
    $ws.Range("L15").Value = -33.333333333333
    $ws.Range("M15").Value = -20
    $ws.Range("F16").Value = 17
    $ws.Range("G16").Value = 13
    $ws.Range("H16").Value = 30.76923076923
    $ws.Range("I16").Value = 122
    $ws.Range("J16").Value = 111
    $ws.Range("K16").Value = 9.909909909909
    $ws.Range("L16").Value = 6.086956521739
    $ws.Range("M16").Value = -48.953974895397
    $ws.Range("N16").Value = -89.756507136859
    $ws.Range("D17").Value = 2
    $ws.Range("E17").Value = 200
    $ws.Range("F17").Value = 33
    $ws.Range("G17").Value = 37
    $ws.Range("H17").Value = -10.81081081081
    $ws.Range("I17").Value = 268
    $ws.Range("J17").Value = 256
    $ws.Range("K17").Value = 4.6875
    $ws.Range("L17").Value = 1.90114068441
    $ws.Range("M17").Value = 0.751879699248
    $ws.Range("N17").Value = -59.759759759759
    $ws.Range("C18").Value = 5
    $ws.Range("E18").Value = 150
    $ws.Range("F18").Value = 14
    $ws.Range("G18").Value = 9
    $ws.Range("H18").Value = 55.555555555555
    $ws.Range("I18").Value = 87
    $ws.Range("J18").Value = 105
    $ws.Range("K18").Value = -17.142857142857
    $ws.Range("L18").Value = -36.95652173913
    $ws.Range("M18").Value = -66.015625
    $ws.Range("N18").Value = -88.274932614555
    $ws.Range("C19").Value = 7
    $ws.Range("D19").Value = 6
    $ws.Range("E19").Value = 16.666666666666
    $ws.Range("F19").Value = 42
    $ws.Range("G19").Value = 35
    $ws.Range("H19").Value = 20
    $ws.Range("I19").Value = 229
    $ws.Range("J19").Value = 229
    $ws.Range("K19").Value = 0
    $ws.Range("L19").Value = -12.260536398467
    $ws.Range("M19").Value = -8.4
    $ws.Range("N19").Value = -54.016064257028
    $ws.Range("C20").Value = 3
    $ws.Range("D20").Value = 5
    $ws.Range("E20").Value = -40
    $ws.Range("F20").Value = 16
    $ws.Range("G20").Value = 11
    $ws.Range("H20").Value = 45.454545454545
    $ws.Range("I20").Value = 68
    $ws.Range("J20").Value = 63
    $ws.Range("K20").Value = 7.936507936507
    $ws.Range("L20").Value = -8.108108108108
    $ws.Range("M20").Value = -9.333333333333
    $ws.Range("N20").Value = -78.881987577639
    $ws.Range("C21").Value = 25
    $ws.Range("D21").Value = 18
    $ws.Range("E21").Value = 38.888888888888
    $ws.Range("F21").Value = 123
    $ws.Range("G21").Value = 107
    $ws.Range("H21").Value = 14.953271028037
    $ws.Range("I21").Value = 792
    $ws.Range("J21").Value = 780
    $ws.Range("K21").Value = 1.538461538461
    $ws.Range("L21").Value = -9.485714285714
    $ws.Range("M21").Value = -28.648648648648
    $ws.Range("N21").Value = -77.345537757437
    $ws.Range("G22").Value = 2
    $ws.Range("H22").Value = 0
    $ws.Range("J22").Value = 9
    $ws.Range("K22").Value = 66.666666666666
    $ws.Range("M22").Value = -40
    $ws.Range("C23").Value = 4
    $ws.Range("D23").Value = 1
    $ws.Range("E23").Value = 300
    $ws.Range("G23").Value = 17
    $ws.Range("H23").Value = -35.294117647058
    $ws.Range("I23").Value = 123
    $ws.Range("J23").Value = 141
    $ws.Range("K23").Value = -12.765957446808
    $ws.Range("L23").Value = -15.753424657534
    $ws.Range("M23").Value = -1.6
    $ws.Range("D24").Value = 23
    $ws.Range("E24").Value = 30.434782608695
    $ws.Range("F24").Value = 128
    $ws.Range("G24").Value = 122
    $ws.Range("H24").Value = 4.918032786885
    $ws.Range("I24").Value = 852
    $ws.Range("J24").Value = 809
    $ws.Range("K24").Value = 5.3152039555
    $ws.Range("L24").Value = -4.591265397536
    $ws.Range("M24").Value = 46.140651801029
    $ws.Range("C25").Value = 13
    $ws.Range("D25").Value = 15
    $ws.Range("E25").Value = -13.333333333333
    $ws.Range("F25").Value = 52
    $ws.Range("G25").Value = 70
    $ws.Range("H25").Value = -25.714285714285
    $ws.Range("I25").Value = 418
    $ws.Range("J25").Value = 351
    $ws.Range("K25").Value = 19.088319088319
    $ws.Range("L25").Value = 21.159420289855
    $ws.Range("C26").Value = 8
    $ws.Range("D26").Value = 10
    $ws.Range("E26").Value = -20
    $ws.Range("F26").Value = 38
    $ws.Range("G26").Value = 42
    $ws.Range("H26").Value = -9.523809523809
    $ws.Range("I26").Value = 346
    $ws.Range("J26").Value = 386
    $ws.Range("K26").Value = -10.362694300518
    $ws.Range("L26").Value = -11.053984575835
    $ws.Range("M26").Value = -39.298245614035
    $ws.Range("L27").Value = -48.148148148148
    $ws.Range("C28").Value = 1
    $ws.Range("I28").Value = 32
    $ws.Range("K28").Value = 6.666666666666
    $ws.Range("L28").Value = 28
    $ws.Range("I29").Value = 13
    $ws.Range("K29").Value = -13.333333333333
    $ws.Range("L29").Value = -7.142857142857
    $ws.Range("M29").Value = -63.888888888888
    $ws.Range("N29").Value = -91.975308641975
    $ws.Range("I30").Value = 12
    $ws.Range("K30").Value = -7.692307692307
    $ws.Range("L30").Value = -7.692307692307
    $ws.Range("M30").Value = -53.846153846153
    $ws.Range("N30").Value = -91.549295774647
    $ws.Range("F31").Value = 3
    $ws.Range("H31").Value = 200
    $ws.Range("I31").Value = 13
    $ws.Range("K31").Value = 116.666666666667
    $ws.Range("L31").Value = 225

# ---------------------------------------------------------------------
# 3. Cells that flip between the text placeholder ("0" / "***.*") and
#    a real number. Excel treats a leading apostrophe as "force text",
#    which is what lets "0"/"***.*" land back in the shared-string
#    table instead of being parsed as a number. After the value is set
#    we restore the correct number style for the column (copy/paste
#    the formatting only) since forcing text temporarily perturbs it.
# ---------------------------------------------------------------------

# Row 22 (Transit): had no week-to-date activity before, now does.
$ws.Range("C22").Value = "'0"
$ws.Range("C14").Copy()
$ws.Range("C22").PasteSpecial(-4122)

$ws.Range("D22").Value = 1
$ws.Range("F14").Copy()
$ws.Range("D22").PasteSpecial(-4122)

$ws.Range("E22").Value = -100
$ws.Range("H14").Copy()
$ws.Range("E22").PasteSpecial(-4122)

# Row 29 (Shooting Vic.): week-to-date count now populated.
$ws.Range("C29").Value = 1
$ws.Range("F14").Copy()
$ws.Range("C29").PasteSpecial(-4122)

# Row 30 (Shooting Inc.): week-to-date count now populated.
$ws.Range("C30").Value = 1
$ws.Range("F14").Copy()
$ws.Range("C30").PasteSpecial(-4122)

# Row 31 (Hate Crimes): week-to-date activity drops back to none.
$ws.Range("D31").Value = "'0"
$ws.Range("C14").Copy()
$ws.Range("D31").PasteSpecial(-4122)

$ws.Range("E31").Value = "'***.*"
$ws.Range("C14").Copy()
$ws.Range("E31").PasteSpecial(-4122)
